# Auto-generated edit script applying TPM recompute diff to Fgf2-Sdc2 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.472738
$ws.Range("H2").Value = 1.418214
$ws.Range("I2").Value = 0.0327564895931267
$ws.Range("J2").Value = 0.03397138804734427
$ws.Range("M2").Value = 1.009860666666667
$ws.Range("N2").Value = 3.029582
$ws.Range("O2").Value = 0.01353413605720072
$ws.Range("P2").Value = 0.01542521070970148
$ws.Range("Q2").Value = 0.4773995118386667
$ws.Range("R2").Value = 4.296595606548
$ws.Range("S2").Value = 0.0004433307869096563
$ws.Range("T2").Value = 0.0005240158187313197
$ws.Range("G3").Value = 0.472738
$ws.Range("H3").Value = 1.418214
$ws.Range("I3").Value = 0.0327564895931267
$ws.Range("J3").Value = 0.03397138804734427
$ws.Range("O3").Value = 0.6185519418990597
$ws.Range("P3").Value = 0.704979911415303
$ws.Range("Q3").Value = 21.81863651003867
$ws.Range("R3").Value = 196.367728590348
$ws.Range("S3").Value = 0.02026159024762486
$ws.Range("T3").Value = 0.02394914613627165
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.472738
$ws.Range("H4").Value = 1.418214
$ws.Range("I4").Value = 0.0327564895931267
$ws.Range("J4").Value = 0.03397138804734427
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009315666666666667
$ws.Range("N4").Value = 0.027947
$ws.Range("O4").Value = 0.0001248484115599408
$ws.Range("P4").Value = 0.000142293017222847
$ws.Range("Q4").Value = 0.004403869628666667
$ws.Range("R4").Value = 0.039634826658
$ws.Range("S4").Value = [double]"4.0895956939816E-06"
$ws.Range("T4").Value = [double]"4.833891304504777E-06"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.472738
$ws.Range("H5").Value = 1.418214
$ws.Range("I5").Value = 0.0327564895931267
$ws.Range("J5").Value = 0.03397138804734427
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 27.4428835
$ws.Range("N5").Value = 54.885767
$ws.Range("O5").Value = 0.3677890736321797
$ws.Range("P5").Value = 0.2794525848577725
$ws.Range("Q5").Value = 12.973293860023
$ws.Range("R5").Value = 77.83976316013801
$ws.Range("S5").Value = 0.0120474789628982
$ws.Range("T5").Value = 0.009493392201036795
$ws.Range("I6").Value = 0.822180234441485
$ws.Range("J6").Value = 0.8526739017519405
$ws.Range("M6").Value = 1.009860666666667
$ws.Range("N6").Value = 3.029582
$ws.Range("O6").Value = 0.01353413605720072
$ws.Range("P6").Value = 0.01542521070970148
$ws.Range("Q6").Value = 11.98261619120889
$ws.Range("R6").Value = 107.84354572088
$ws.Range("S6").Value = 0.01112749915647225
$ws.Range("T6").Value = 0.01315267460118698
$ws.Range("I7").Value = 0.822180234441485
$ws.Range("J7").Value = 0.8526739017519405
$ws.Range("O7").Value = 0.6185519418990597
$ws.Range("P7").Value = 0.704979911415303
$ws.Range("S7").Value = 0.5085611806048047
$ws.Range("T7").Value = 0.6011179717232238
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.822180234441485
$ws.Range("J8").Value = 0.8526739017519405
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.009315666666666667
$ws.Range("N8").Value = 0.027947
$ws.Range("O8").Value = 0.0001248484115599408
$ws.Range("P8").Value = 0.000142293017222847
$ws.Range("Q8").Value = 0.1105360986088889
$ws.Range("R8").Value = 0.99482488748
$ws.Range("S8").Value = 0.0001026478962859991
$ws.Range("T8").Value = 0.000121329542187461
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.822180234441485
$ws.Range("J9").Value = 0.8526739017519405
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.4428835
$ws.Range("N9").Value = 54.885767
$ws.Range("O9").Value = 0.3677890736321797
$ws.Range("P9").Value = 0.2794525848577725
$ws.Range("Q9").Value = 325.6266443627134
$ws.Range("R9").Value = 1953.75986617628
$ws.Range("S9").Value = 0.3023889067839221
$ws.Range("T9").Value = 0.2382819258853421
$ws.Range("G10").Value = 0.37892
$ws.Range("H10").Value = 1.13676
$ws.Range("I10").Value = 0.02625574638939025
$ws.Range("J10").Value = 0.02722954016579943
$ws.Range("M10").Value = 1.009860666666667
$ws.Range("N10").Value = 3.029582
$ws.Range("O10").Value = 0.01353413605720072
$ws.Range("P10").Value = 0.01542521070970148
$ws.Range("Q10").Value = 0.3826564038133333
$ws.Range("R10").Value = 3.44390763432
$ws.Range("S10").Value = 0.0003553488439173642
$ws.Range("T10").Value = 0.000420021394585736
$ws.Range("G11").Value = 0.37892
$ws.Range("H11").Value = 1.13676
$ws.Range("I11").Value = 0.02625574638939025
$ws.Range("J11").Value = 0.02722954016579943
$ws.Range("O11").Value = 0.6185519418990597
$ws.Range("P11").Value = 0.704979911415303
$ws.Range("Q11").Value = 17.48858299181333
$ws.Range("R11").Value = 157.39724692632
$ws.Range("S11").Value = 0.01624054291516656
$ws.Range("T11").Value = 0.01919627881396472
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 0.37892
$ws.Range("H12").Value = 1.13676
$ws.Range("I12").Value = 0.02625574638939025
$ws.Range("J12").Value = 0.02722954016579943
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.009315666666666667
$ws.Range("N12").Value = 0.027947
$ws.Range("O12").Value = 0.0001248484115599408
$ws.Range("P12").Value = 0.000142293017222847
$ws.Range("Q12").Value = 0.003529892413333333
$ws.Range("R12").Value = 0.03176903172
$ws.Range("S12").Value = [double]"3.277988231036023E-06"
$ws.Range("T12").Value = [double]"3.874573427782302E-06"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 0.37892
$ws.Range("H13").Value = 1.13676
$ws.Range("I13").Value = 0.02625574638939025
$ws.Range("J13").Value = 0.02722954016579943
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 27.4428835
$ws.Range("N13").Value = 54.885767
$ws.Range("O13").Value = 0.3677890736321797
$ws.Range("P13").Value = 0.2794525848577725
$ws.Range("Q13").Value = 10.39865741582
$ws.Range("R13").Value = 62.39194449492
$ws.Range("S13").Value = 0.009656576642075285
$ws.Range("T13").Value = 0.007609365383821191
$ws.Range("G14").Value = 1.548357
$ws.Range("H14").Value = 3.096714
$ws.Range("I14").Value = 0.1072872076222874
$ws.Range("J14").Value = 0.0741775733180209
$ws.Range("M14").Value = 1.009860666666667
$ws.Range("N14").Value = 3.029582
$ws.Range("O14").Value = 0.01353413605720072
$ws.Range("P14").Value = 0.01542521070970148
$ws.Range("Q14").Value = 1.563624832258
$ws.Range("R14").Value = 9.381748993548001
$ws.Range("S14").Value = 0.001452039665157179
$ws.Range("T14").Value = 0.001144204698364803
$ws.Range("G15").Value = 1.548357
$ws.Range("H15").Value = 3.096714
$ws.Range("I15").Value = 0.1072872076222874
$ws.Range("J15").Value = 0.0741775733180209
$ws.Range("O15").Value = 0.6185519418990597
$ws.Range("P15").Value = 0.704979911415303
$ws.Range("Q15").Value = 71.462498404558
$ws.Range("R15").Value = 428.774990427348
$ws.Range("S15").Value = 0.06636271061569345
$ws.Range("T15").Value = 0.05229369906674051
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 1.548357
$ws.Range("H16").Value = 3.096714
$ws.Range("I16").Value = 0.1072872076222874
$ws.Range("J16").Value = 0.0741775733180209
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.009315666666666667
$ws.Range("N16").Value = 0.027947
$ws.Range("O16").Value = 0.0001248484115599408
$ws.Range("P16").Value = 0.000142293017222847
$ws.Range("Q16").Value = 0.014423977693
$ws.Range("R16").Value = 0.08654386615800001
$ws.Range("S16").Value = [double]"1.339463745234415E-05"
$ws.Range("T16").Value = [double]"1.055495071769014E-05"
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 1.548357
$ws.Range("H17").Value = 3.096714
$ws.Range("I17").Value = 0.1072872076222874
$ws.Range("J17").Value = 0.0741775733180209
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 27.4428835
$ws.Range("N17").Value = 54.885767
$ws.Range("O17").Value = 0.3677890736321797
$ws.Range("P17").Value = 0.2794525848577725
$ws.Range("Q17").Value = 42.49138076740951
$ws.Range("R17").Value = 169.965523069638
$ws.Range("S17").Value = 0.03945906270398439
$ws.Range("T17").Value = 0.02072911460219788
$ws.Range("G18").Value = 0.16626
$ws.Range("H18").Value = 0.49878
$ws.Range("I18").Value = 0.01152032195371061
$ws.Range("J18").Value = 0.01194759671689489
$ws.Range("M18").Value = 1.009860666666667
$ws.Range("N18").Value = 3.029582
$ws.Range("O18").Value = 0.01353413605720072
$ws.Range("P18").Value = 0.01542521070970148
$ws.Range("Q18").Value = 0.16789943444
$ws.Range("R18").Value = 1.51109490996
$ws.Range("S18").Value = 0.0001559176047442758
$ws.Range("T18").Value = 0.0001842941968326414
$ws.Range("G19").Value = 0.16626
$ws.Range("H19").Value = 0.49878
$ws.Range("I19").Value = 0.01152032195371061
$ws.Range("J19").Value = 0.01194759671689489
$ws.Range("O19").Value = 0.6185519418990597
$ws.Range("P19").Value = 0.704979911415303
$ws.Range("Q19").Value = 7.673524248439999
$ws.Range("R19").Value = 69.06171823595999
$ws.Range("S19").Value = 0.007125917515770065
$ws.Range("T19").Value = 0.008422815675102327
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 0.16626
$ws.Range("H20").Value = 0.49878
$ws.Range("I20").Value = 0.01152032195371061
$ws.Range("J20").Value = 0.01194759671689489
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.009315666666666667
$ws.Range("N20").Value = 0.027947
$ws.Range("O20").Value = 0.0001248484115599408
$ws.Range("P20").Value = 0.000142293017222847
$ws.Range("Q20").Value = 0.00154882274
$ws.Range("R20").Value = 0.01393940466
$ws.Range("S20").Value = [double]"1.438293896579883E-06"
$ws.Range("T20").Value = [double]"1.700059585408755E-06"
$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = 0.16626
$ws.Range("H21").Value = 0.49878
$ws.Range("I21").Value = 0.01152032195371061
$ws.Range("J21").Value = 0.01194759671689489
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 27.4428835
$ws.Range("N21").Value = 54.885767
$ws.Range("O21").Value = 0.3677890736321797
$ws.Range("P21").Value = 0.2794525848577725
$ws.Range("Q21").Value = 4.56265381071
$ws.Range("R21").Value = 27.37592286426
$ws.Range("S21").Value = 0.004237048539299687
$ws.Range("T21").Value = 0.003338786785374515
